$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2 changes from "ECs" to "MuSCs"; other numeric values updated
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.3821933333333334
$ws.Range("H2").Value = 1.14658
$ws.Range("I2").Value = 0.1675988106211496
$ws.Range("J2").Value = 0.1675988106211496
$ws.Range("M2").Value = 0.000484
$ws.Range("N2").Value = 0.001452
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0001849815733333334
$ws.Range("R2").Value = 0.00166483416
$ws.Range("S2").Value = 0.1675988106211496
$ws.Range("T2").Value = 0.1675988106211496

# Row 3: A3 changes from "ECs" to "FAPs"; other numeric values updated
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 1.440755333333333
$ws.Range("H3").Value = 4.322266
$ws.Range("I3").Value = 0.6317977295855797
$ws.Range("J3").Value = 0.6317977295855797
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.0006973255813333334
$ws.Range("R3").Value = 0.006275930232
$ws.Range("S3").Value = 0.6317977295855797
$ws.Range("T3").Value = 0.6317977295855797

# Row 4: A4 changes from "FAPs" to "MuSCs"; D4 changes from "ECs" to "MuSCs"
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.4574573333333333
$ws.Range("H4").Value = 1.372372
$ws.Range("I4").Value = 0.2006034597932707
$ws.Range("J4").Value = 0.2006034597932707
$ws.Range("M4").Value = 0.000484
$ws.Range("N4").Value = 0.001452
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.0002214093493333333
$ws.Range("R4").Value = 0.001992684144
$ws.Range("S4").Value = 0.2006034597932707
$ws.Range("T4").Value = 0.2006034597932707

# Remove rows 5, 6, 7 (the old extra rows for MuSCs/FAPs combos)
$ws.Range("A5:T7").EntireRow.Delete()
